$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1819.6
$ws.Range("I19").Value = 2700
$ws.Range("J19").Value = 1232.6666
$ws.Range("K19").Value = 2700
$ws.Range("L19").Value = 1232.6666
$ws.Range("M19").Value = -2525
$ws.Range("N19").Value = -1582.6666
$ws.Range("H100").Value = 1600.375
$ws.Range("I100").Value = 1600.375
$ws.Range("K100").Value = 1600.375
$ws.Range("M100").Value = -1059.375
$ws.Range("H135").Value = 1000
$ws.Range("I135").Value = 1000
$ws.Range("K135").Value = 9000
$ws.Range("M135").Value = -6465

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5549.6875
$ws.Range("I32").Value = 4155.769
$ws.Range("K32").Value = 4155.769
$ws.Range("M32").Value = -3868.769
$ws.Range("H74").Value = 824.3333
$ws.Range("I74").Value = 824.3333
$ws.Range("K74").Value = 824.3333
$ws.Range("M74").Value = 49.66669999999999
$ws.Range("H77").Value = 824.3333
$ws.Range("I77").Value = 824.3333
$ws.Range("K77").Value = 4121.6665
$ws.Range("M77").Value = 246.3334999999997
$ws.Range("H80").Value = 24605
$ws.Range("J80").Value = 28110
$ws.Range("L80").Value = 28110
$ws.Range("N80").Value = -30106
$ws.Range("H83").Value = 24605
$ws.Range("J83").Value = 28110
$ws.Range("L83").Value = 84330
$ws.Range("N83").Value = -94314
$ws.Range("H110").Value = 2000
$ws.Range("I110").Value = 2000
$ws.Range("K110").Value = 2000
$ws.Range("M110").Value = 45

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4405.8335
$ws.Range("I86").Value = 4721.7144
$ws.Range("K86").Value = 4721.7144
$ws.Range("M86").Value = -3598.7144
$ws.Range("H89").Value = 4405.8335
$ws.Range("I89").Value = 4721.7144
$ws.Range("K89").Value = 23608.572
$ws.Range("M89").Value = -17992.572
$ws.Range("H99").Value = 3317
$ws.Range("I99").Value = 2707.4
$ws.Range("J99").Value = 4333
$ws.Range("K99").Value = 2707.4
$ws.Range("L99").Value = 4333
$ws.Range("M99").Value = -1209.4
$ws.Range("N99").Value = -7329

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1959.5
$ws.Range("I31").Value = 1621.6666
$ws.Range("K31").Value = 1621.6666
$ws.Range("M31").Value = -1326.6666
$ws.Range("H34").Value = 1959.5
$ws.Range("I34").Value = 1621.6666
$ws.Range("K34").Value = 1621.6666
$ws.Range("M34").Value = -1419.6666
$ws.Range("H105").Value = 1927.7142
$ws.Range("I105").Value = 1698.8
$ws.Range("K105").Value = 1698.8
$ws.Range("M105").Value = 48.20000000000005
$ws.Range("H122").Value = 3612.125
$ws.Range("I122").Value = 3483.6667
$ws.Range("K122").Value = 10451.0001
$ws.Range("M122").Value = -8001.000100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 223.45454
$ws.Range("J7").Value = 390.66666
$ws.Range("L7").Value = 1171.99998
$ws.Range("N7").Value = -1395.99998
$ws.Range("H23").Value = 222.28572
$ws.Range("J23").Value = 271.2
$ws.Range("L23").Value = 813.5999999999999
$ws.Range("N23").Value = -1283.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10375
$ws.Range("I70").Value = 10375
$ws.Range("K70").Value = 10375
$ws.Range("M70").Value = -10105
$ws.Range("H73").Value = 10375
$ws.Range("I73").Value = 10375
$ws.Range("K73").Value = 10375
$ws.Range("M73").Value = -9439
$ws.Range("H80").Value = 11745.125
$ws.Range("I80").Value = 12985
$ws.Range("J80").Value = 11001.2
$ws.Range("K80").Value = 12985
$ws.Range("L80").Value = 11001.2
$ws.Range("M80").Value = -11987
$ws.Range("N80").Value = -12997.2
$ws.Range("H83").Value = 11745.125
$ws.Range("I83").Value = 12985
$ws.Range("J83").Value = 11001.2
$ws.Range("K83").Value = 64925
$ws.Range("L83").Value = 55006
$ws.Range("M83").Value = -59933
$ws.Range("N83").Value = -64990
$ws.Range("H113").Value = 1314.3334
$ws.Range("I113").Value = 1321.75
$ws.Range("K113").Value = 1321.75
$ws.Range("M113").Value = 848.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1212.909
$ws.Range("I16").Value = 1253.2858
$ws.Range("J16").Value = 1142.25
$ws.Range("K16").Value = 1253.2858
$ws.Range("L16").Value = 1142.25
$ws.Range("M16").Value = -1083.2858
$ws.Range("N16").Value = -1482.25
$ws.Range("H68").Value = 1100
$ws.Range("I68").Value = 1100
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1100
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -351
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 1100
$ws.Range("I71").Value = 1100
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 5500
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -1756
$ws.Range("N71").ClearContents()
$ws.Range("H82").Value = 863.1667
$ws.Range("I82").Value = 919.75
$ws.Range("J82").Value = 750
$ws.Range("K82").Value = 919.75
$ws.Range("L82").Value = 750
$ws.Range("M82").Value = -558.75
$ws.Range("N82").Value = -1472
$ws.Range("H85").Value = 863.1667
$ws.Range("I85").Value = 919.75
$ws.Range("J85").Value = 750
$ws.Range("K85").Value = 919.75
$ws.Range("L85").Value = 750
$ws.Range("M85").Value = 328.25
$ws.Range("N85").Value = -3246
$ws.Range("H100").Value = 4062.5
$ws.Range("I100").Value = 4083.3333
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 4083.3333
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -3542.3333
$ws.Range("N100").Value = -5082

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H107").Value = 594.3333
$ws.Range("I107").Value = 578.5714
$ws.Range("K107").Value = 1735.7142
$ws.Range("M107").Value = 184.2857999999999
$ws.Range("H113").Value = 466
$ws.Range("I113").Value = 466
$ws.Range("K113").Value = 1398
$ws.Range("M113").Value = 772
$ws.Range("H116").Value = 50000
$ws.Range("J116").Value = 50000
$ws.Range("L116").Value = 50000
$ws.Range("N116").Value = -59178
